# Update the arithmetic problems in the document to match the new
# generated output (commit c986bee).
#
# Each "XX÷Y=" expression is unique within the document, so a simple
# Find/Replace (no wildcards, match whole expression) for each pair is
# sufficient and safe against accidental partial matches.

$d = $word.ActiveDocument

$replacements = @(
    @("74÷2=", "21÷3="),
    @("74÷5=", "57÷7="),
    @("46÷4=", "59÷6="),
    @("35÷4=", "31÷3="),
    @("86÷8=", "95÷4="),
    @("94÷8=", "31÷3="),
    @("71÷9=", "72÷9="),
    @("29÷6=", "84÷8="),
    @("89÷5=", "94÷3="),
    @("61÷2=", "59÷7="),
    @("53÷8=", "78÷8="),
    @("83÷8=", "39÷8="),
    @("87÷4=", "71÷7="),
    @("15÷6=", "23÷2="),
    @("67÷3=", "36÷3="),
    @("38÷8=", "10÷3="),
    @("99÷3=", "66÷2="),
    @("45÷9=", "53÷3="),
    @("26÷7=", "26÷2="),
    @("88÷6=", "26÷4="),
    @("45÷5=", "78÷4="),
    @("41÷6=", "98÷3="),
    @("50÷9=", "36÷5="),
    @("94÷4=", "71÷2="),
    @("66÷4=", "44÷2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]

    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true,
                         1, $false, $new, 2)
}

$d.Save()
